$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table of exposure sites (replaces the old Cheltenham/Moorabbin rows).
$data = @(
    @("Cheltenham",   "Bodero Southland Shopping Centre, 1239 Nepean Hwy",             "22/12/20 6.45pm - 7pm",     "Case shopped in store", "new"),
    @("Cheltenham",   "Chemist Warehouse Cheltenham, 326/330 Charman Rd",               "03/01/21, 3.30pm - 3.45pm", "Case shopped in store", "new"),
    @("Cheltenham",   "Coles, Westfield Southland",                                     "22/12/20 11:50am-12:10pm",  "Case shopped in store", "new"),
    @("Cheltenham",   "Kmart Southland Shopping Centre, 1239 Nepean Highway",           "22/12/20 6.30pm - 6.45pm",  "Case shopped in store", "new"),
    @("Cheltenham",   "Kmart Southland Shopping Centre, 1239 Nepean Highway",           "28/12/20 2.30pm-3pm",       "Case shopped in store", "new"),
    @("Cheltenham",   "Specsavers, 1004-1005 Westfield Southland",                      "22/12/20 11:00am-1145am",   "Case shopped in store", "new"),
    @("Dandenong",    "Kmart - Clayton Plaza, 2107 Dandenong Rd",                       "30/12/20 7pm - 7.30pm",     "Case shopped at store", "new"),
    @("Dandenong",    "Woolworths - Clayton Plaza, 2107 Dandenong Rd",                  "30/12/20 7.30pm - 745pm",   "Case shopped at store", "new"),
    @("Frankston",    "TK Maxx Frankston, 10 Shannon Street, Bayside Shopping Centre",  "31/12/20 2pm - 3pm",        "Case shopped at store", "new"),
    @("Mount Martha", "Mount Martha Fine Foods, 34 Lochiel Ave",                        "31/12/20 3pm - 3.15pm",     "Takeaway coffee",       "new"),
    @("Springvale",   "IKEA Springvale, 917 Princes Hwy",                               "29/12/20 4pm - 6pm",        "Case shopped at store and dined at cafe", "new")
)

# Clear out the old data rows (2-6) completely before laying down the new table.
$ws.Rows.Item(2).Resize(5,1).EntireRow.ClearContents()

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Resize columns to (best-effort) match the content-driven widths the workbook
# ships with after the new, wider table is written.
$ws.Columns.Item(1).ColumnWidth = 11.5
$ws.Columns.Item(2).ColumnWidth = 51.166666666666664
$ws.Columns.Item(3).ColumnWidth = 22.666666666666668
$ws.Columns.Item(4).ColumnWidth = 31.500000000000004
